# Automatische test-sync: 2025-09-02 23:08:50
#
# This script reproduces the addition of a new log entry ("Bestelling M6
# bouten" placed by planning@testbedrijf123.nl) to the "Logs" sheet, and
# updates the "Dashboard" sheet's summary table (row order/counts) plus the
# conditional-formatting ranges on "Logs" that need to grow to include the
# new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new row (row 4) to the "Logs" sheet
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Bestelling M6 bouten"
$logs.Range("B4").Value = '"Testbedrijf 123 B.V." <planning@testbedrijf123.nl>'
$logs.Range("D4").Value = "Inkoop / Bestellingen"
$logs.Range("F4").Value = "2025-09-02 23:08:13"
$logs.Range("G4").Value = "Nee"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting ranges on "Logs" so they cover the
#    new row (2:3 -> 2:4) for columns D, G, H, I, J
# ---------------------------------------------------------------------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "3")
    $newRange = $logs.Range($col + "2:" + $col + "4")
    $conds = $oldRange.FormatConditions
    for ($i = 1; $i -le $conds.Count; $i++) {
        $conds.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Update the "Dashboard" summary table: the two category rows swap
#    order and the "Inkoop / Bestellingen" count increases to 2
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Inkoop / Bestellingen"
$dash.Range("B2").Value = 2
$dash.Range("A3").Value = "Klacht / Probleem"
$dash.Range("B3").Value = 1
